# Generate Report for Handback
# - Marks the zh-cn / de-de handback status as "Handed back: in sync with en-US"
# - Fills in the "Latest Target File" (source .md, hyperlinked) and
#   "Latest Handback File" (generated .xlf) columns for both source docs
# - Stamps the "Latest Handback DateTime" with the real handback timestamps
#   (zh-cn and de-de were handed back at slightly different times)

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$file1Name = "438699ce-5cb5-4e86-822d-2b503fedfca4.md"
$file1Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ce5e95b7211bb252ad6cd1e81e6c43b8cff8716/e2e/438699ce-5cb5-4e86-822d-2b503fedfca4.md"

$file2Name = "dc90e2a8-89e2-4731-b5ab-66a7f86b0c93.md"
$file2Url  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6ce5e95b7211bb252ad6cd1e81e6c43b8cff8716/e2e/dc90e2a8-89e2-4731-b5ab-66a7f86b0c93.md"

# ----- Overview sheet: update the per-language status cells -----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $statusText
$ov.Range("F2").Value = $statusText
$ov.Range("E3").Value = $statusText
$ov.Range("F3").Value = $statusText
$ov.Columns.Item(5).ColumnWidth = 29.166666666666668
$ov.Columns.Item(6).ColumnWidth = 29.166666666666668

# ----- zh-cn sheet -----
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

$zh.Hyperlinks.Add($zh.Range("I2"), $file1Url, "", "", $file1Name)
$zh.Range("J2").Value = "438699ce-5cb5-4e86-822d-2b503fedfca4.7bdecfbf4e8cacf10ad7deeb9e1485049a6645be.zh-cn.xlf"
$zh.Range("K2").Value = "2016-08-27 16:59:40"

$zh.Hyperlinks.Add($zh.Range("I3"), $file2Url, "", "", $file2Name)
$zh.Range("J3").Value = "dc90e2a8-89e2-4731-b5ab-66a7f86b0c93.6faa6f6cf1fce11535fd9e4732dfea1382638e2b.zh-cn.xlf"
$zh.Range("K3").Value = "2016-08-27 16:59:40"

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ----- de-de sheet -----
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

$de.Hyperlinks.Add($de.Range("I2"), $file1Url, "", "", $file1Name)
$de.Range("J2").Value = "438699ce-5cb5-4e86-822d-2b503fedfca4.7bdecfbf4e8cacf10ad7deeb9e1485049a6645be.de-de.xlf"
$de.Range("K2").Value = "2016-08-27 16:59:47"

$de.Hyperlinks.Add($de.Range("I3"), $file2Url, "", "", $file2Name)
$de.Range("J3").Value = "dc90e2a8-89e2-4731-b5ab-66a7f86b0c93.6faa6f6cf1fce11535fd9e4732dfea1382638e2b.de-de.xlf"
$de.Range("K3").Value = "2016-08-27 16:59:47"

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664
